$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D2: append the "cloud production" navigation step
$ws.Range("D2").Value = "输入用户名`"longchangkun`"`n输入密码`"Lck123456`"`n勾选我已阅读并同意复选框`n点击`"登录`"按钮`n点击左侧导航栏的“统计”下面的 “云制作”"

# Update E2: change URL target to /desktop
$ws.Range("E2").Value = "跳转成功到页面，https://task-pre.renderbus.com/desktop"

# Add new row 13 content
$ws.Range("D13").Value = "输入用户名`"longchangkun`"`n输入密码`"Lck123456`"`n勾选我已阅读并同意复选框`n点击`"登录`"按钮"
$ws.Range("D13").WrapText = $true
$ws.Range("E13").Value = "跳转成功到页面，https://task-pre.renderbus.com/"
$ws.Rows(13).RowHeight = 68

# Clear old rows 17-21 (the exploded step rows no longer needed)
$ws.Range("D17:D21").Clear()

# Update the active cell selection
[void]$ws.Range("D6").Select()
